$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("C7").Value = 7
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 10000
$ws.Range("G7").Value = 0.2
$ws.Range("H7").Value = 26
$ws.Range("I7").Value = "naïve test release same as before"

# Row 8
$ws.Range("C8").Value = 7
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 10000
$ws.Range("G8").Value = 0.2
$ws.Range("H8").Value = 108.2
$ws.Range("I8").Value = "scattered"

# Update selection to match target state
$ws.Range("H9").Select()
